$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("values")
$rng = $ws1.Range("AA1:AA1048576")
$fc = $rng.FormatConditions.Add(1, 3, "0.637")
$fc.Interior.Color = 65535
Write-Output "added"
